# Apply the edit described by the diff:
# For data rows 2-8:
#   Column B (position) -> -1
#   Column C (type)      -> "date"
#   Column U (csim)       -> 0   (only actually changes on rows 2,5,6,7 which were 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $ws.Range("B$row").Value = -1
    $ws.Range("C$row").Value = "date"
    $ws.Range("U$row").Value = 0
}
